$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: beteckning, förändrad-datum, area
$ws.Range("A2").Value = "A 36523-2022"
$ws.Range("C2").Value = 46065
$ws.Range("G2").Value = 0.2

# Row 3: beteckning, förändrad-datum, area
$ws.Range("A3").Value = "A 36578-2022"
$ws.Range("C3").Value = 46065
$ws.Range("G3").Value = 0.3
